# finish Appendices 13th April
# Insert a new title row above the data table, merge it into four blocks
# (B:E, F:H, I:K, L:N) and label the first block "Bryum_X1", centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing header/data rows down by one to make room for the title row
$ws.Rows.Item(1).Insert()

# Re-create the four merged blocks across the new title row
$ws.Range("B1:E1").Merge()
$ws.Range("F1:H1").Merge()
$ws.Range("I1:K1").Merge()
$ws.Range("L1:N1").Merge()

# Center the text across the whole title row
$ws.Range("B1:N1").HorizontalAlignment = -4108

# Label the first merged block
$ws.Range("B1").Value = "Bryum_X1"

# Restore the selection/view state captured in the saved workbook
$ws.Range("F1:N16").Select()

Write-Host "Applied Bryum_X1 title row edit"
